# Add new weekly ranking sheet for 2025-11-26, appended after the last existing sheet.
$wb = $excel.ActiveWorkbook

$titles = @(
    'ブルーロック',
    'みいちゃんと山田さん',
    'ガチアクタ',
    '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！',
    'ギルティサークル',
    '東京卍リベンジャーズ',
    '島耕作',
    '薫る花は凛と咲く',
    'FAIRY TAIL 100 YEARS QUEST',
    '十字架のろくにん',
    '愛妻の裏アカ',
    '転生したら第七王子だったので、気ままに魔術を極めます',
    'ハードワーカー中田',
    'WIND BREAKER',
    '南海トラフ巨大地震',
    '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～',
    '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す',
    '君が僕らを悪魔と呼んだ頃',
    '魔女と傭兵',
    '異世界ウォーキング',
    '蒼く染めろ',
    'アルキメデスの大戦',
    'K-9~警視庁公安部公安第9課異能対策係~',
    'さわらないで小手指くん',
    'ドラハチ',
    'ひゃくえむ。',
    '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～',
    'グラぱらっ！',
    'おやすみ ふみさん',
    '屋根の下のアルテミス',
    '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜',
    'ハナバス　苔石花江のバスケ論',
    'せいぶつ部の田辺くん',
    'イレギュラーズ',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    'となりの黒川さん',
    '幼馴染とはラブコメにならない',
    '食糧人類-Starving Anonymous-',
    'ナキナギ',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    '阿武ノーマル',
    '不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～',
    'ジュミドロ',
    'アオバノバスケ',
    '黄昏町プリズナーズ',
    'いじめるヤバイ奴',
    '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！',
    'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。',
    '可愛いだけじゃない式守さん',
    '春くらり',
    'デッドアカウント',
    '最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    'デスティニーラバーズ',
    '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜',
    '時々ボソッとロシア語でデレる隣のアーリャさん',
    '五輪の女神さま 〜なでしこ寮のメダルごはん〜',
    'ストーカー行為がバレて人生終了男',
    'ともだちづくり',
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
    '黒猫と魔女の教室',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    'インフェクション',
    'MYS',
    '劣等人の魔剣使い　スキルボードを駆使して最強に至る',
    'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜',
    'それがメイドのカンナです',
    '四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～',
    '東京卍リベンジャーズ～場地圭介からの手紙～',
    '恋ニ非ズ',
    '母という呪縛 娘という牢獄',
    'リスナーに騙されてダンジョンの最下層から脱出RTAすることになった',
    '不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～',
    '降り積もれ孤独な死よ',
    '勇者と呼ばれた後に　―そして無双男は家族を創る―',
    '魁の花巫女',
    '我間乱 ―修羅―',
    '復讐の教科書',
    '東大リベンジャーズ',
    'なれの果ての僕ら',
    'ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～',
    'この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～',
    '触手魔術師の成り上がり',
    '転生したらスライムだった件',
    'ハプスブルク家の華麗なる受難',
    '剣帝学院の魔眼賢者',
    '金田一少年の事件簿外伝 犯人たちの事件簿',
    'お願い、脱がシて。',
    '君が監督！',
    'ヒロインは絶望しました。',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    'GALAXIAS',
    '人間消失',
    '田んぼで拾った女騎士、田舎で俺の嫁だと思われている',
    'ザ・ファブル',
    '中華一番！極',
    '絶望集落',
    '彼女、お借りします',
    'イジらないで、長瀞さん',
    'お嬢様の僕'
)

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "magapoke_2025-11-26"

# Header row
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"

$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data rows: rank number in column A, title in column B
for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}

# Restore original active sheet/selection (sheet order unaffected by active-tab state)
$wb.Worksheets.Item(1).Activate()
